$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# "Tobias" -> "Tobias Rempel" for rows 29-45 (B column), and fill in his
# Matrikelnummer in column C for the same rows.
$ws.Range("B29:B45").Value = "Tobias Rempel"
$ws.Range("C29:C45").Value = 7080879

# "Andy" -> "Andy Kruder" for rows 46-49 and 52-54 (B column).
$ws.Range("B46:B49").Value = "Andy Kruder"
$ws.Range("B52:B54").Value = "Andy Kruder"

# Update the view state to match where the author left the selection.
$ws.Range("C47").Select()
